$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "mitigation sources" worksheet right after "mitigation"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "mitigation sources"

# Row 13: US_NY mitigation date, sourced from the Imperial College report
[void]$ws2.Hyperlinks.Add($ws2.Range("A1"), "https://www.imperial.ac.uk/mrc-global-infectious-disease-analysis/covid-19/report-13-europe-npi-impact/")

[void]$ws1.Range("B2").Copy($ws1.Range("B13"))
$ws1.Range("A13").Value = "US_NY"
$ws1.Range("B13").Value = "3/12/2020"

# Row 14: US_LA mitigation date, sourced from the CDC report
[void]$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.cdc.gov/mmwr/volumes/69/wr/mm6915e2.htm")

[void]$ws1.Range("B2").Copy($ws1.Range("B14"))
$ws1.Range("A14").Value = "US_LA"
$ws1.Range("B14").Value = "3/13/2020"

# Row 15: US_WA mitigation date
[void]$ws1.Range("B2").Copy($ws1.Range("B15"))
$ws1.Range("A15").Value = "US_WA"
$ws1.Range("B15").Value = "3/15/2020"

# Row 16: US_CA mitigation date
[void]$ws1.Range("B2").Copy($ws1.Range("B16"))
$ws1.Range("A16").Value = "US_CA"
$ws1.Range("B16").Value = "3/11/2020"

[void]$ws1.Range("C16").Select()
[void]$ws2.Range("A2").Select()
[void]$ws1.Select()
